$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for rows 2-5 from 2023-10-22 to 2023-10-25
$ws.Range("C2:C5").Value = 45224
